# Add "Source Links" sheet, update the G3 note on "Data Boosting Notes",
# and wrap/resize that row to fit the longer text.

$wb = $excel.ActiveWorkbook
$wsNotes = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Update the "Notes:" cell (G3) on the existing sheet with the
#    expanded note, enable wrap text and grow the row to fit it.
# ---------------------------------------------------------------------
$wsNotes.Range("G3").Value = "Dumped everything into trash - couldn't generalize. This iteration added 25,000 images. Will reevaluate the weighting issue after retraining in next iteration."
$wsNotes.Range("G3").WrapText = $true
$wsNotes.Rows.Item(3).RowHeight = 48.75

# ---------------------------------------------------------------------
# 2) Insert a new "Sheet1" worksheet right after "Data Boosting Notes"
#    and fill it in with the source links that were used for this
#    iteration.
# ---------------------------------------------------------------------
$wsLinks = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsNotes)
$wsLinks.Name = "Sheet1"
$wsLinks.Columns.Item(1).ColumnWidth = 24.11

# Header
$wsLinks.Range("A1").Value = "Source Links:"
$wsLinks.Range("A1").Font.Underline = $true
$wsLinks.Range("A1").Font.Color = 16711680

# Plain hyperlink rows
$wsLinks.Range("A2").Value = "https://www.kaggle.com/moezabid/bottles-and-cans"
$wsLinks.Hyperlinks.Add($wsLinks.Range("A2"), "https://www.kaggle.com/moezabid/bottles-and-cans") | Out-Null
$wsLinks.Range("A2").Font.Underline = $true
$wsLinks.Range("A2").Font.Color = 16711680

$wsLinks.Range("A3").Value = "https://medium.com/@ringlayer/cardboard-box-detection-using-retinanet-keras-5d4f331d9d15"
$wsLinks.Hyperlinks.Add($wsLinks.Range("A3"), "https://medium.com/@ringlayer/cardboard-box-detection-using-retinanet-keras-5d4f331d9d15") | Out-Null
$wsLinks.Range("A3").Font.Underline = $true
$wsLinks.Range("A3").Font.Color = 16711680

$wsLinks.Range("A4").Value = "https://www.kaggle.com/techsash/waste-classification-data"
$wsLinks.Hyperlinks.Add($wsLinks.Range("A4"), "https://www.kaggle.com/techsash/waste-classification-data") | Out-Null
$wsLinks.Range("A4").Font.Underline = $true
$wsLinks.Range("A4").Font.Color = 16711680

$wsLinks.Range("A5").Value = "https://www.kaggle.com/asdasdasasdas/garbage-classification"
$wsLinks.Hyperlinks.Add($wsLinks.Range("A5"), "https://www.kaggle.com/asdasdasasdas/garbage-classification") | Out-Null
$wsLinks.Range("A5").Font.Underline = $true
$wsLinks.Range("A5").Font.Color = 16711680

# Mixed-format row: plain text prefix followed by a hyperlink-styled URL
$prefixText = "DuckDuck Go Scraper - "
$urlText = "https://colab.research.google.com/github/joedockrill/image-scraper/blob/master/ImageScraper.ipynb"
$wsLinks.Range("A6").Value = $prefixText + $urlText
$wsLinks.Hyperlinks.Add($wsLinks.Range("A6"), $urlText) | Out-Null
$wsLinks.Range("A6").Font.Underline = $true
$wsLinks.Range("A6").Font.Color = 16711680

$rt = $wsLinks.Range("A6").Characters($prefixText.Length + 1, $urlText.Length)
$rt.Font.Underline = $true
$rt.Font.Color = 0xCC5511

Write-Host "Done"
